$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "Apple"
$ws.Range("B3").Value = "'0.29"
$ws.Range("C3").Value = "'0.64"
$ws.Range("D3").Value = "'13.11"
$ws.Range("E3").Value = "'62.64"
$ws.Range("A4").Value = "Beans"
$ws.Range("B4").Value = "'3.71"
$ws.Range("C4").Value = "'0.6"
$ws.Range("D4").Value = "'2.85"
$ws.Range("E4").Value = "'33.46"
$ws.Range("A5").Value = "Butter"
$ws.Range("B5").Value = "'0.85"
$ws.Range("C5").Value = "'81.11"
$ws.Range("D5").Value = "'0.06"
$ws.Range("E5").Value = "'717"
$ws.Range("A6").Value = "Capsicum"
$ws.Range("B6").Value = "'1.11"
$ws.Range("C6").Value = "'0.34"
$ws.Range("D6").Value = "'1.84"
$ws.Range("E6").Value = "'16.25"
$ws.Range("A7").Value = "Carrot"
$ws.Range("B7").Value = "'1.04"
$ws.Range("C7").Value = "'0.47"
$ws.Range("D7").Value = "'6.71"
$ws.Range("E7").Value = "'38.24"
$ws.Range("A8").Value = "Chicken"
$ws.Range("B8").Value = "'21.81"
$ws.Range("C8").Value = "'9"
$ws.Range("D8").Value = "'0"
$ws.Range("E8").Value = "'168.26"
$ws.Range("A9").Value = "Coriander"
$ws.Range("B9").Value = "'3.52"
$ws.Range("C9").Value = "'0.7"
$ws.Range("D9").Value = "'1.93"
$ws.Range("E9").Value = "'31.07"
$ws.Range("A10").Value = "Corn"
$ws.Range("B10").Value = "'2.69"
$ws.Range("C10").Value = "'1.33"
$ws.Range("D10").Value = "'11.66"
$ws.Range("E10").Value = "'73.14"
$ws.Range("A11").Value = "Corn flour"
$ws.Range("B11").Value = "'1.2"
$ws.Range("C11").Value = "'0.2"
$ws.Range("D11").Value = "'88.4"
$ws.Range("E11").Value = "'360"
$ws.Range("A12").Value = "Cucumber"
$ws.Range("B12").Value = "'0.71"
$ws.Range("C12").Value = "'0.16"
$ws.Range("D12").Value = "'3.48"
$ws.Range("E12").Value = "'19.6"
$ws.Range("A13").Value = "Curd"
$ws.Range("B13").Value = "'3.1"
$ws.Range("C13").Value = "'4"
$ws.Range("D13").Value = "'3"
$ws.Range("E13").Value = "'60"
$ws.Range("A14").Value = "Custard Powder"
$ws.Range("B14").Value = "'1.9"
$ws.Range("C14").Value = "'1.5"
$ws.Range("D14").Value = "'84"
$ws.Range("E14").Value = "'357"
$ws.Range("A15").Value = "Dalia"
$ws.Range("B15").Value = "'10.84"
$ws.Range("C15").Value = "'1.45"
$ws.Range("D15").Value = "'69.06"
$ws.Range("E15").Value = "'341.78"
$ws.Range("A16").Value = "Egg (boiled)"
$ws.Range("B16").Value = "'12.37"
$ws.Range("C16").Value = "'0.26"
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'52.58"
$ws.Range("A17").Value = "Egg (raw)"
$ws.Range("B17").Value = "'10.84"
$ws.Range("C17").Value = "'0.06"
$ws.Range("D17").Value = "'0"
$ws.Range("E17").Value = "'44.69"
$ws.Range("A18").Value = "Gourd"
$ws.Range("B18").Value = "'0.49"
$ws.Range("C18").Value = "'0.13"
$ws.Range("D18").Value = "'2.25"
$ws.Range("E18").Value = "'12.91"
$ws.Range("A19").Value = "Guava"
$ws.Range("B19").Value = "'1.44"
$ws.Range("C19").Value = "'0.32"
$ws.Range("D19").Value = "'5.13"
$ws.Range("E19").Value = "'32.27"
$ws.Range("A20").Value = "Lentil"
$ws.Range("B20").Value = "'22.87"
$ws.Range("C20").Value = "'0.61"
$ws.Range("D20").Value = "'47.91"
$ws.Range("E20").Value = "'297.8"
$ws.Range("A21").Value = "Lettuce"
$ws.Range("B21").Value = "'1.54"
$ws.Range("C21").Value = "'0.27"
$ws.Range("D21").Value = "'3.01"
$ws.Range("E21").Value = "'21.75"
$ws.Range("A22").Value = "Mango"
$ws.Range("B22").Value = "'0.46"
$ws.Range("C22").Value = "'0.54"
$ws.Range("D22").Value = "'9.03"
$ws.Range("E22").Value = "'44.69"
$ws.Range("A23").Value = "Milk"
$ws.Range("B23").Value = "'3.26"
$ws.Range("C23").Value = "'4.48"
$ws.Range("D23").Value = "'4.94"
$ws.Range("E23").Value = "'72.9"
$ws.Range("A24").Value = "Oats"
$ws.Range("B24").Value = "'13.6"
$ws.Range("C24").Value = "'7.6"
$ws.Range("D24").Value = "'62.8"
$ws.Range("E24").Value = "'374"
$ws.Range("A25").Value = "Oil"
$ws.Range("B25").Value = "'0"
$ws.Range("C25").Value = "'100"
$ws.Range("D25").Value = "'0"
$ws.Range("E25").Value = "'900"
$ws.Range("A26").Value = "Onion"
$ws.Range("B26").Value = "'1.5"
$ws.Range("C26").Value = "'0"
$ws.Range("D26").Value = "'9.56"
$ws.Range("E26").Value = "'48.04"
$ws.Range("A27").Value = "Paneer"
$ws.Range("B27").Value = "'18.86"
$ws.Range("C27").Value = "'14.78"
$ws.Range("D27").Value = "'12.41"
$ws.Range("E27").Value = "'257.89"
$ws.Range("A28").Value = "Papaya (raw)"
$ws.Range("B28").Value = "'0.5"
$ws.Range("C28").Value = "'0.23"
$ws.Range("D28").Value = "'4.4"
$ws.Range("E28").Value = "'23.92"
$ws.Range("A29").Value = "Papaya (ripe)"
$ws.Range("B29").Value = "'0"
$ws.Range("C29").Value = "'0"
$ws.Range("D29").Value = "'4.61"
$ws.Range("E29").Value = "'23.9"
$ws.Range("A30").Value = "Pear"
$ws.Range("B30").Value = "'0"
$ws.Range("C30").Value = "'0"
$ws.Range("D30").Value = "'8.09"
$ws.Range("E30").Value = "'37.52"
$ws.Range("A31").Value = "Peas"
$ws.Range("B31").Value = "'7.25"
$ws.Range("C31").Value = "'0"
$ws.Range("D31").Value = "'11.88"
$ws.Range("E31").Value = "'81.26"
$ws.Range("A32").Value = "Pineapple"
$ws.Range("B32").Value = "'0"
$ws.Range("C32").Value = "'0"
$ws.Range("D32").Value = "'9.42"
$ws.Range("E32").Value = "'43.02"
$ws.Range("A33").Value = "Potato"
$ws.Range("B33").Value = "'1.54"
$ws.Range("C33").Value = "'0"
$ws.Range("D33").Value = "'14.89"
$ws.Range("E33").Value = "'69.79"
$ws.Range("A34").Value = "Puffed rice"
$ws.Range("B34").Value = "'7.47"
$ws.Range("C34").Value = "'1.62"
$ws.Range("D34").Value = "'77.68"
$ws.Range("E34").Value = "'361.85"
$ws.Range("A35").Value = "Rice"
$ws.Range("B35").Value = "'7.94"
$ws.Range("C35").Value = "'0.52"
$ws.Range("D35").Value = "'78.24"
$ws.Range("E35").Value = "'356.36"
$ws.Range("A36").Value = "Rice flakes"
$ws.Range("B36").Value = "'7.44"
$ws.Range("C36").Value = "'1.14"
$ws.Range("D36").Value = "'76.75"
$ws.Range("E36").Value = "'353.73"
$ws.Range("A37").Value = "Semolina"
$ws.Range("B37").Value = "'11.38"
$ws.Range("C37").Value = "'0"
$ws.Range("D37").Value = "'68.43"
$ws.Range("E37").Value = "'333.65"
$ws.Range("A38").Value = "Skimmed Milk"
$ws.Range("B38").Value = "'3.37"
$ws.Range("C38").Value = "'0.08"
$ws.Range("D38").Value = "'4.89"
$ws.Range("E38").Value = "'34"
$ws.Range("A39").Value = "Spinach"
$ws.Range("B39").Value = "'2.14"
$ws.Range("C39").Value = "'0.64"
$ws.Range("D39").Value = "'2.05"
$ws.Range("E39").Value = "'24.38"
$ws.Range("A40").Value = "Sugar"
$ws.Range("B40").Value = "'0.1"
$ws.Range("C40").Value = "'0"
$ws.Range("D40").Value = "'99.4"
$ws.Range("E40").Value = "'398"
$ws.Range("A41").Value = "Sweet potato"
$ws.Range("B41").Value = "'1.27"
$ws.Range("C41").Value = "'0"
$ws.Range("D41").Value = "'23.93"
$ws.Range("E41").Value = "'108.03"
$ws.Range("A42").Value = "Tomato"
$ws.Range("B42").Value = "'0.9"
$ws.Range("C42").Value = "'0.47"
$ws.Range("D42").Value = "'2.71"
$ws.Range("E42").Value = "'19.6"
$ws.Range("A43").Value = "Vermicelli"
$ws.Range("B43").Value = "'9.7"
$ws.Range("C43").Value = "'0"
$ws.Range("D43").Value = "'70.39"
$ws.Range("E43").Value = "'332.7"
$ws.Range("A44").Value = "Wheat flour"
$ws.Range("B44").Value = "'29.2"
$ws.Range("C44").Value = "'7.4"
$ws.Range("D44").Value = "'53.3"
$ws.Range("E44").Value = "'397"
$ws.Range("A1:E44").Style = "Normal"